# Apply "player_data" updates described in the commit:
# "use snaking algorithm rather than recursive + update player data"
#
# The algorithm change itself isn't representable in the worksheet (it
# lives in whatever external tool generated the assignments), but the
# resulting player-data corrections in the sheet are:
#
#   Row 14 (Edie Hetling):   remove G14 ("C") and H14 ("Y")
#   Row 16 (Kathie Moen):    E16/F16  "C"  -> "C-"
#   Row 20 (Lisa Tsaur):     J20 gains "Sinndy DeJesus", existing value
#                            ("Aly Werth") shifts to K20
#   Row 27 (Nick Bryant):    F27      "B-" -> "B"
#   Row 28 (Chris Butts):    E28      "C+" -> "B-"
#   Row 31 (Roe Hendrick):   E31 "B" -> "B-", F31 "B-" -> "C+"
#   Row 44 (Scott Riggot):   E44 "B+" -> "B", remove I44 ("Y")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: drop the extra Offense/Defense rating columns G14 and H14
$ws.Range("G14").ClearContents()
$ws.Range("H14").ClearContents()

# Row 16: rating "C" -> "C-"
$ws.Range("E16").Value = "C-"
$ws.Range("F16").Value = "C-"

# Row 20: add a second conflict (Sinndy DeJesus) ahead of the existing one
$ws.Range("K20").Value = $ws.Range("J20").Value2
$ws.Range("J20").Value = "Sinndy DeJesus"

# Row 27: rating "B-" -> "B"
$ws.Range("F27").Value = "B"

# Row 28: rating "C+" -> "B-"
$ws.Range("E28").Value = "B-"

# Row 31: rating "B" -> "B-", "B-" -> "C+"
$ws.Range("E31").Value = "B-"
$ws.Range("F31").Value = "C+"

# Row 44: rating "B+" -> "B"; drop the extra Setting rating column I44
$ws.Range("E44").Value = "B"
$ws.Range("I44").ClearContents()
